$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 96 (shifts "phone", "password", "user_update" rows
# of the Users table down by one, and shifts the whole Roles table down by one row)
$ws.Rows(96).Insert()

# Copy the formatting of the (now shifted) "phone" row's Type cell onto the new
# row's Type cell, since the plain insert does not carry the correct fill/border
# formatting down from the table body.
$ws.Range("B97").Copy()
$ws.Range("B96").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row with the "address" column definition
$ws.Range("A96").Value = "address"
$ws.Range("B96").Value = "String"
$ws.Range("C96").Value = 255
$ws.Range("D96").Formula = '=IF(B96="String", CONCATENATE("$table->",LOWER(B96),"(''",A96,"'', ",C96,");"), IF(B96="Integer", CONCATENATE("$table->",LOWER(B96),"(''",A96,"'')->unsigned()->default(0);"), IF(B96="Text", CONCATENATE("$table->",LOWER(B96),"(''",A96,"'');"), IF(B96="Date", CONCATENATE("$table->","timestamp","(''",A96,"'');")) )))'

# Grow the "Users" table (Table25468) to include the new row
$loUsers = $ws.ListObjects("Table25468")
$loUsers.Resize($ws.Range("A89:D99"))

# The "Roles" table (Table254687) has moved down by one row
$loRoles = $ws.ListObjects("Table254687")
$loRoles.Resize($ws.Range("A104:D106"))

# Update view/selection state to mirror what a user would see after this edit
$win = $excel.ActiveWindow
$win.ScrollRow = 85
$win.ScrollColumn = 1
$ws.Range("D96").Select()

# Page setup tweaks (paper size set to A4/Letter as in the saved file)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
